$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: swap UART_TX/UART_RX between I3/J3 ---
$ws.Range("I3").Value = "UART_RX"
$ws.Range("J3").Value = "UART_TX"

# --- New row 22: ESP32 Type A ---
$ws.Range("B22").Value = "ESP32 Type A"
$ws.Range("C22").Value = "GND"
$ws.Range("D22").Value = "5V"
$ws.Range("E22").Value = "PIR"
$ws.Range("F22").Value = "I2C_Data"
$ws.Range("G22").Value = "I2C_Clock"
$ws.Range("H22").Value = "3V3"

# --- New row 26: 6P / colours ---
$ws.Range("B26").Value = "6P"
$ws.Range("C26").Value = "Green"
$ws.Range("D26").Value = "Orange"
$ws.Range("E26").Value = "Yellow"
$ws.Range("F26").Value = "White"
$ws.Range("G26").Value = "Red"
$ws.Range("H26").Value = "Black"

# --- New row 27: I2D/I2C mapping ---
$ws.Range("C27").Value = "I2D"
$ws.Range("D27").Value = "I2C"
$ws.Range("E27").Value = "PIR"
$ws.Range("F27").Value = "3V3"
$ws.Range("G27").Value = "5V"
$ws.Range("H27").Value = "GND"

# --- Row 3: add the note to K3 (added after the other new strings so it is
#     appended last in the shared-strings table) ---
$ws.Range("K3").Value = "UART_TX will be dominant colour, as its important to identify output from esp pin (note this is transmit out of esp, thus RX from peripherials connect to this)"

# --- View adjustments ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("K4").Select()
